$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-10-14 Tuesday" "2025-10-15 Wednesday"

Replace-Text "25×25=625" "44×80=3520"
Replace-Text "15×99=1485" "62×45=2790"
Replace-Text "57×70=3990" "95×88=8360"
Replace-Text "52×76=3952" "83×87=7221"
Replace-Text "47×77=3619" "31×54=1674"

Replace-Text "91×47=4277" "60×69=4140"
Replace-Text "78×80=6240" "64×28=1792"
Replace-Text "72×63=4536" "11×87=957"
Replace-Text "35×66=2310" "42×91=3822"
Replace-Text "66×46=3036" "11×37=407"

Replace-Text "96×94=9024" "50×18=900"
Replace-Text "79×29=2291" "62×67=4154"
Replace-Text "54×36=1944" "79×18=1422"
Replace-Text "97×80=7760" "74×26=1924"
Replace-Text "99×43=4257" "61×75=4575"

Replace-Text "11×49=539" "77×41=3157"
Replace-Text "94×17=1598" "86×59=5074"
Replace-Text "27×72=1944" "77×62=4774"
Replace-Text "72×61=4392" "42×75=3150"
Replace-Text "84×88=7392" "93×21=1953"

Replace-Text "33×80=2640" "42×77=3234"
Replace-Text "87×31=2697" "74×31=2294"
Replace-Text "72×84=6048" "17×66=1122"
Replace-Text "77×56=4312" "15×97=1455"
Replace-Text "37×38=1406" "93×30=2790"
